$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per diff
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.993.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.57%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.47%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.49%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4304"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07250"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.121.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +19.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8683"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.44%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.421"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.09%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.623"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06952"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.88%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.84%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008901"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.11%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.035.03"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.18%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.367.58"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +18.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.887"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.79%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.228"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.68%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +11.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.83"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08956"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.189"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7482"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.433"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.815"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.007"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.48%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05239"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01928"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.77%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5111"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.65%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1658"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.741"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.501"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.358"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.67%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.42"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.62%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.650"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.23%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4571"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.46%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06223"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.836"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.25%  "
